$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 41.428665
$ws.Range("H2").Value = 124.285995
$ws.Range("I2").Value = 0.06969137269740189
$ws.Range("J2").Value = 0.06969137269740189
$ws.Range("M2").Value = 173.5452066666667
$ws.Range("N2").Value = 520.63562
$ws.Range("O2").Value = 0.6098887991422922
$ws.Range("P2").Value = 0.6098887991422922
$ws.Range("Q2").Value = 7189.7462293491
$ws.Range("R2").Value = 64707.7160641419
$ws.Range("S2").Value = 0.04250398760499637
$ws.Range("T2").Value = 0.04250398760499637
$ws.Range("G3").Value = 41.428665
$ws.Range("H3").Value = 124.285995
$ws.Range("I3").Value = 0.06969137269740189
$ws.Range("J3").Value = 0.06969137269740189
$ws.Range("O3").Value = 0.003264284357140855
$ws.Range("P3").Value = 0.003264284357140855
$ws.Range("Q3").Value = 38.48140215279
$ws.Range("R3").Value = 346.33261937511
$ws.Range("S3").Value = 0.0002274924577238022
$ws.Range("T3").Value = 0.0002274924577238022
$ws.Range("G4").Value = 41.428665
$ws.Range("H4").Value = 124.285995
$ws.Range("I4").Value = 0.06969137269740189
$ws.Range("J4").Value = 0.06969137269740189
$ws.Range("M4").Value = 54.64271666666667
$ws.Range("N4").Value = 163.92815
$ws.Range("O4").Value = 0.192030546333187
$ws.Range("P4").Value = 0.192030546333187
$ws.Range("Q4").Value = 2263.77480347325
$ws.Range("R4").Value = 20373.97323125925
$ws.Range("S4").Value = 0.01338287237379184
$ws.Range("T4").Value = 0.01338287237379184
$ws.Range("G5").Value = 41.428665
$ws.Range("H5").Value = 124.285995
$ws.Range("I5").Value = 0.06969137269740189
$ws.Range("J5").Value = 0.06969137269740189
$ws.Range("M5").Value = 1.069012
$ws.Range("N5").Value = 3.207036
$ws.Range("O5").Value = 0.00375682196858928
$ws.Range("P5").Value = 0.00375682196858928
$ws.Range("Q5").Value = 44.28774002898
$ws.Range("R5").Value = 398.58966026082
$ws.Range("S5").Value = 0.0002618180799707425
$ws.Range("T5").Value = 0.0002618180799707426
$ws.Range("G6").Value = 41.428665
$ws.Range("H6").Value = 124.285995
$ws.Range("I6").Value = 0.06969137269740189
$ws.Range("J6").Value = 0.06969137269740189
$ws.Range("M6").Value = 54.36641700000001
$ws.Range("N6").Value = 163.099251
$ws.Range("O6").Value = 0.1910595481987908
$ws.Range("P6").Value = 0.1910595481987908
$ws.Range("Q6").Value = 2252.328077143305
$ws.Range("R6").Value = 20270.95269428974
$ws.Range("S6").Value = 0.01331520218091915
$ws.Range("T6").Value = 0.01331520218091915
$ws.Range("H7").Value = 510.696747
$ws.Range("I7").Value = 0.2863649869040173
$ws.Range("J7").Value = 0.2863649869040173
$ws.Range("M7").Value = 173.5452066666667
$ws.Range("N7").Value = 520.63562
$ws.Range("O7").Value = 0.6098887991422922
$ws.Range("P7").Value = 0.6098887991422922
$ws.Range("Q7").Value = 29542.99083403646
$ws.Range("R7").Value = 265886.9175063281
$ws.Range("S7").Value = 0.1746507979792893
$ws.Range("T7").Value = 0.1746507979792894
$ws.Range("H8").Value = 510.696747
$ws.Range("I8").Value = 0.2863649869040173
$ws.Range("J8").Value = 0.2863649869040173
$ws.Range("O8").Value = 0.003264284357140855
$ws.Range("P8").Value = 0.003264284357140855
$ws.Range("S8").Value = 0.0009347767471836294
$ws.Range("T8").Value = 0.0009347767471836295
$ws.Range("H9").Value = 510.696747
$ws.Range("I9").Value = 0.2863649869040173
$ws.Range("J9").Value = 0.2863649869040173
$ws.Range("M9").Value = 54.64271666666667
$ws.Range("N9").Value = 163.92815
$ws.Range("O9").Value = 0.192030546333187
$ws.Range("P9").Value = 0.192030546333187
$ws.Range("Q9").Value = 9301.952549636451
$ws.Range("R9").Value = 83717.57294672805
$ws.Range("S9").Value = 0.05499082488587439
$ws.Range("T9").Value = 0.0549908248858744
$ws.Range("H10").Value = 510.696747
$ws.Range("I10").Value = 0.2863649869040173
$ws.Range("J10").Value = 0.2863649869040173
$ws.Range("M10").Value = 1.069012
$ws.Range("N10").Value = 3.207036
$ws.Range("O10").Value = 0.00375682196858928
$ws.Range("P10").Value = 0.00375682196858928
$ws.Range("Q10").Value = 181.980316967988
$ws.Range("R10").Value = 1637.822852711892
$ws.Range("S10").Value = 0.001075822273835794
$ws.Range("T10").Value = 0.001075822273835794
$ws.Range("H11").Value = 510.696747
$ws.Range("I11").Value = 0.2863649869040173
$ws.Range("J11").Value = 0.2863649869040173
$ws.Range("M11").Value = 54.36641700000001
$ws.Range("N11").Value = 163.099251
$ws.Range("O11").Value = 0.1910595481987908
$ws.Range("P11").Value = 0.1910595481987908
$ws.Range("Q11").Value = 9254.917435981833
$ws.Range("R11").Value = 83294.2569238365
$ws.Range("S11").Value = 0.05471276501783418
$ws.Range("T11").Value = 0.05471276501783419
$ws.Range("G12").Value = 244.5761666666666
$ws.Range("H12").Value = 733.7284999999999
$ws.Range("I12").Value = 0.4114264551867299
$ws.Range("J12").Value = 0.41142645518673
$ws.Range("M12").Value = 173.5452066666667
$ws.Range("N12").Value = 520.63562
$ws.Range("O12").Value = 0.6098887991422922
$ws.Range("P12").Value = 0.6098887991422922
$ws.Range("Q12").Value = 42445.02138990777
$ws.Range("R12").Value = 382005.19250917
$ws.Range("S12").Value = 0.2509243866892048
$ws.Range("T12").Value = 0.2509243866892049
$ws.Range("G13").Value = 244.5761666666666
$ws.Range("H13").Value = 733.7284999999999
$ws.Range("I13").Value = 0.4114264551867299
$ws.Range("J13").Value = 0.41142645518673
$ws.Range("O13").Value = 0.003264284357140855
$ws.Range("P13").Value = 0.003264284357140855
$ws.Range("Q13").Value = 227.1768551192222
$ws.Range("R13").Value = 2044.591696073
$ws.Range("S13").Value = 0.001343012941779955
$ws.Range("T13").Value = 0.001343012941779956
$ws.Range("G14").Value = 244.5761666666666
$ws.Range("H14").Value = 733.7284999999999
$ws.Range("I14").Value = 0.4114264551867299
$ws.Range("J14").Value = 0.41142645518673
$ws.Range("M14").Value = 54.64271666666667
$ws.Range("N14").Value = 163.92815
$ws.Range("O14").Value = 0.192030546333187
$ws.Range("P14").Value = 0.192030546333187
$ws.Range("Q14").Value = 13364.30617858611
$ws.Range("R14").Value = 120278.755607275
$ws.Range("S14").Value = 0.07900644696543424
$ws.Range("T14").Value = 0.07900644696543425
$ws.Range("G15").Value = 244.5761666666666
$ws.Range("H15").Value = 733.7284999999999
$ws.Range("I15").Value = 0.4114264551867299
$ws.Range("J15").Value = 0.41142645518673
$ws.Range("M15").Value = 1.069012
$ws.Range("N15").Value = 3.207036
$ws.Range("O15").Value = 0.00375682196858928
$ws.Range("P15").Value = 0.00375682196858928
$ws.Range("Q15").Value = 261.4548570806667
$ws.Range("R15").Value = 2353.093713726
$ws.Range("S15").Value = 0.00154565594530432
$ws.Range("T15").Value = 0.00154565594530432
$ws.Range("G16").Value = 244.5761666666666
$ws.Range("H16").Value = 733.7284999999999
$ws.Range("I16").Value = 0.4114264551867299
$ws.Range("J16").Value = 0.41142645518673
$ws.Range("M16").Value = 54.36641700000001
$ws.Range("N16").Value = 163.099251
$ws.Range("O16").Value = 0.1910595481987908
$ws.Range("P16").Value = 0.1910595481987908
$ws.Range("Q16").Value = 13296.7298652615
$ws.Range("R16").Value = 119670.5687873535
$ws.Range("S16").Value = 0.07860695264500665
$ws.Range("T16").Value = 0.07860695264500667
$ws.Range("G17").Value = 24.173247
$ws.Range("H17").Value = 72.51974100000001
$ws.Range("I17").Value = 0.04066427836821081
$ws.Range("J17").Value = 0.04066427836821081
$ws.Range("M17").Value = 173.5452066666667
$ws.Range("N17").Value = 520.63562
$ws.Range("O17").Value = 0.6098887991422922
$ws.Range("P17").Value = 0.6098887991422922
$ws.Range("Q17").Value = 4195.151146419381
$ws.Range("R17").Value = 37756.36031777443
$ws.Range("S17").Value = 0.02480068790197598
$ws.Range("T17").Value = 0.02480068790197598
$ws.Range("G18").Value = 24.173247
$ws.Range("H18").Value = 72.51974100000001
$ws.Range("I18").Value = 0.04066427836821081
$ws.Range("J18").Value = 0.04066427836821081
$ws.Range("O18").Value = 0.003264284357140855
$ws.Range("P18").Value = 0.003264284357140855
$ws.Range("Q18").Value = 22.453546092922
$ws.Range("R18").Value = 202.081914836298
$ws.Range("S18").Value = 0.0001327397677717718
$ws.Range("T18").Value = 0.0001327397677717718
$ws.Range("G19").Value = 24.173247
$ws.Range("H19").Value = 72.51974100000001
$ws.Range("I19").Value = 0.04066427836821081
$ws.Range("J19").Value = 0.04066427836821081
$ws.Range("M19").Value = 54.64271666666667
$ws.Range("N19").Value = 163.92815
$ws.Range("O19").Value = 0.192030546333187
$ws.Range("P19").Value = 0.192030546333187
$ws.Range("Q19").Value = 1320.89188673435
$ws.Range("R19").Value = 11888.02698060915
$ws.Range("S19").Value = 0.007808783591292322
$ws.Range("T19").Value = 0.007808783591292322
$ws.Range("G20").Value = 24.173247
$ws.Range("H20").Value = 72.51974100000001
$ws.Range("I20").Value = 0.04066427836821081
$ws.Range("J20").Value = 0.04066427836821081
$ws.Range("M20").Value = 1.069012
$ws.Range("N20").Value = 3.207036
$ws.Range("O20").Value = 0.00375682196858928
$ws.Range("P20").Value = 0.00375682196858928
$ws.Range("Q20").Value = 25.841491121964
$ws.Range("R20").Value = 232.5734200976761
$ws.Range("S20").Value = 0.0001527684543105242
$ws.Range("T20").Value = 0.0001527684543105242
$ws.Range("G21").Value = 24.173247
$ws.Range("H21").Value = 72.51974100000001
$ws.Range("I21").Value = 0.04066427836821081
$ws.Range("J21").Value = 0.04066427836821081
$ws.Range("M21").Value = 54.36641700000001
$ws.Range("N21").Value = 163.099251
$ws.Range("O21").Value = 0.1910595481987908
$ws.Range("P21").Value = 0.1910595481987908
$ws.Range("Q21").Value = 1314.212826645999
$ws.Range("R21").Value = 11827.91543981399
$ws.Range("S21").Value = 0.007769298652860218
$ws.Range("T21").Value = 0.007769298652860218
$ws.Range("G22").Value = 114.0486906666667
$ws.Range("H22").Value = 342.146072
$ws.Range("I22").Value = 0.19185290684364
$ws.Range("J22").Value = 0.19185290684364
$ws.Range("M22").Value = 173.5452066666667
$ws.Range("N22").Value = 520.63562
$ws.Range("O22").Value = 0.6098887991422922
$ws.Range("P22").Value = 0.6098887991422922
$ws.Range("Q22").Value = 19792.60359180941
$ws.Range("R22").Value = 178133.4323262846
$ws.Range("S22").Value = 0.1170089389668256
$ws.Range("T22").Value = 0.1170089389668256
$ws.Range("G23").Value = 114.0486906666667
$ws.Range("H23").Value = 342.146072
$ws.Range("I23").Value = 0.19185290684364
$ws.Range("J23").Value = 0.19185290684364
$ws.Range("O23").Value = 0.003264284357140855
$ws.Range("P23").Value = 0.003264284357140855
$ws.Range("Q23").Value = 105.9351907801796
$ws.Range("R23").Value = 953.416717021616
$ws.Range("S23").Value = 0.0006262624426816956
$ws.Range("T23").Value = 0.0006262624426816956
$ws.Range("G24").Value = 114.0486906666667
$ws.Range("H24").Value = 342.146072
$ws.Range("I24").Value = 0.19185290684364
$ws.Range("J24").Value = 0.19185290684364
$ws.Range("M24").Value = 54.64271666666667
$ws.Range("N24").Value = 163.92815
$ws.Range("O24").Value = 0.192030546333187
$ws.Range("P24").Value = 0.192030546333187
$ws.Range("Q24").Value = 6231.930290302978
$ws.Range("R24").Value = 56087.37261272681
$ws.Range("S24").Value = 0.03684161851679423
$ws.Range("T24").Value = 0.03684161851679423
$ws.Range("G25").Value = 114.0486906666667
$ws.Range("H25").Value = 342.146072
$ws.Range("I25").Value = 0.19185290684364
$ws.Range("J25").Value = 0.19185290684364
$ws.Range("M25").Value = 1.069012
$ws.Range("N25").Value = 3.207036
$ws.Range("O25").Value = 0.00375682196858928
$ws.Range("P25").Value = 0.00375682196858928
$ws.Range("Q25").Value = 121.9194189069547
$ws.Range("R25").Value = 1097.274770162592
$ws.Range("S25").Value = 0.0007207572151678993
$ws.Range("T25").Value = 0.0007207572151678994
$ws.Range("G26").Value = 114.0486906666667
$ws.Range("H26").Value = 342.146072
$ws.Range("I26").Value = 0.19185290684364
$ws.Range("J26").Value = 0.19185290684364
$ws.Range("M26").Value = 54.36641700000001
$ws.Range("N26").Value = 163.099251
$ws.Range("O26").Value = 0.1910595481987908
$ws.Range("P26").Value = 0.1910595481987908
$ws.Range("Q26").Value = 6200.418675088009
$ws.Range("R26").Value = 55803.76807579208
$ws.Range("S26").Value = 0.03665532970217054
$ws.Range("T26").Value = 0.03665532970217054
